$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (F1:I1): new standardized-variable labels ---
$ws.Range("F1").Value = "Ypadr"
$ws.Range("G1").Value = "X1padr"
$ws.Range("H1").Value = "X2padr"
$ws.Range("I1").Value = "X3padr"

# Format header cells like the other header cells (bold, centered)
$hdr = $ws.Range("F1:I1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108

# --- Standardized-value formulas (F2:I19) ---
$ws.Range("F2:F19").Formula = "=(B2-AVERAGE(`$B`$2:`$B`$19))/STDEV.S(`$B`$2:`$B`$19)"
$ws.Range("G2:G19").Formula = "=(C2-AVERAGE(`$C`$2:`$C`$19))/STDEV.S(`$C`$2:`$C`$19)"
$ws.Range("H2:H19").Formula = "=(D2-AVERAGE(`$D`$2:`$D`$19))/STDEV.S(`$D`$2:`$D`$19)"
$ws.Range("I2:I19").Formula = "=(E2-AVERAGE(`$E`$2:`$E`$19))/STDEV.S(`$E`$2:`$E`$19)"

# Last data row (19) gets the same bottom border styling as A19:E19
$bottom = $ws.Range("F19:I19").Borders.Item(9)
$bottom.LineStyle = 1
$bottom.Weight = 2

# Keep the selection where the author left it
$ws.Range("K9").Select()
